$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Dodatočné informácie:" column (old column E) is removed entirely; this
# shifts the old "Očakávaný výsledok:" column (old F) left into the new
# column E, and drops the now-unused shared string.
$ws.Columns("E").Delete()

# The remaining "Testované dáta" column (D) is narrowed and the results
# column (new E) is widened to fit its content after the layout change.
$ws.Columns("D").ColumnWidth = 15.333333333333334
$ws.Columns("E").ColumnWidth = 57.5

# Re-fit the row heights of the data rows now that the text wraps into the
# new column widths.
$ws.Rows(2).RowHeight = 19.5
$ws.Rows(4).RowHeight = 15.75
$ws.Rows(7).RowHeight = 29.25
$ws.Rows(8).RowHeight = 40.5
$ws.Rows(9).RowHeight = 123
$ws.Rows(10).RowHeight = 40.5
$ws.Rows(11).RowHeight = 47.25
$ws.Rows(12).RowHeight = 49.5

# Page is now printed in landscape (xlLandscape = 2).
$ws.PageSetup.Orientation = 2

# Update the remembered selection.
$ws.Range("G6").Select()
